$d = $word.ActiveDocument

# 1) Remove the "Meta description" paragraph (currently paragraph #2,
#    right after the title heading).
$metaPara = $d.Paragraphs.Item(2)
[void]$metaPara.Range.Delete()

# 2) Replace the closing "image prompt" paragraph with two paragraphs:
#    - a bold "Play 7s Wild Free..." paragraph (re-using the removed
#      meta description's bold headline text)
#    - the former meta-description body text, in italics (matching the
#      formatting of the paragraph it replaces)
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$targetRange = $lastPara.Range
[void]$targetRange.MoveEnd(1, 1)

$xmlFrag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play 7s Wild Free: Simple Gameplay with Wild Symbols and Free Spins</w:t></w:r></w:p>
<w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of 7s Wild, a simple slot game with wild symbols and free spins, and play it for free. Discover pros, cons, and similar games.</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

[void]$targetRange.InsertXML($xmlFrag)
